$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a numeric-looking string that must remain text
# (matches the source data which stores these as text, not numbers).
$textCells = @(
    "D5",
    "D6",
    "D7",
    "D9",
    "D10",
    "D12",
    "D15",
    "D16",
    "D20",
    "D21",
    "D22",
    "D23",
    "D24",
    "D26",
    "D28",
    "D29",
    "D30",
    "D31",
    "D32",
    "D34",
    "D35",
    "D36",
    "D37",
    "D38",
    "D39",
    "D40",
    "D41",
    "D45",
    "D46",
    "D48",
    "D49",
    "D51",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$updates = [ordered]@{
    "D2" = "44.400.62"
    "E2" = "  +0.76%  "
    "D3" = "2.244.41"
    "E3" = "  +0.08%  "
    "E4" = "  +0.31%  "
    "D5" = "307.84"
    "E5" = "  +0.66%  "
    "D6" = "94.32"
    "E6" = "  -2.38%  "
    "D7" = "0.571"
    "E7" = "  -0.25%  "
    "E8" = "  +0.20%  "
    "D9" = "0.525"
    "E9" = "  -0.18%  "
    "D10" = "34.94"
    "E10" = "  +0.48%  "
    "E11" = "  +0.30%  "
    "D12" = "7.20"
    "E12" = "  +0.24%  "
    "E13" = "  +0.91%  "
    "D14" = "2.358.39"
    "E14" = "  +5.27%  "
    "D15" = "0.839"
    "E15" = "  +2.06%  "
    "D16" = "13.68"
    "E16" = "  +0.40%  "
    "D17" = "44.090.11"
    "E17" = "  +0.47%  "
    "D18" = "0.0₃0966"
    "E18" = "  -0.02%  "
    "E19" = "  -0.45%  "
    "D20" = "6.41"
    "E20" = "  +2.94%  "
    "D21" = "65.82"
    "E21" = "  +1.59%  "
    "D22" = "3.02"
    "E22" = "  +3.39%  "
    "D23" = "237.33"
    "E23" = "  -0.71%  "
    "D24" = "2.00"
    "E24" = "  +3.45%  "
    "E25" = "  -0.21%  "
    "D26" = "38.43"
    "E26" = "  +5.53%  "
    "E27" = "  +4.65%  "
    "D28" = "9.85"
    "E28" = "  -0.78%  "
    "D29" = "5.95"
    "E29" = "  -1.35%  "
    "D30" = "20.06"
    "E30" = "  +0.27%  "
    "D31" = "153.95"
    "E31" = "  +0.17%  "
    "D32" = "0.0799"
    "E32" = "  -0.95%  "
    "E33" = "  -0.21%  "
    "D34" = "3.10"
    "E34" = "  -7.64%  "
    "B35" = "Kaspa"
    "C35" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "D35" = "0.109"
    "E35" = "  +2.99%  "
    "B36" = "Stellar"
    "C36" = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
    "D36" = "0.120"
    "E36" = "  +1.08%  "
    "D37" = "1.80"
    "E37" = "  +1.93%  "
    "D38" = "3.48"
    "E38" = "  +5.35%  "
    "D39" = "14.66"
    "E39" = "  -1.24%  "
    "D40" = "3.81"
    "E40" = "  +0.63%  "
    "D41" = "0.0304"
    "E41" = "  +0.59%  "
    "E42" = "  +0.36%  "
    "D43" = "1.744.09"
    "E43" = "  -0.29%  "
    "E44" = "  +2.66%  "
    "D45" = "80.53"
    "E45" = "  -6.10%  "
    "D46" = "99.88"
    "E46" = "  -0.54%  "
    "E47" = "  -2.85%  "
    "D48" = "70.95"
    "E48" = "  +3.18%  "
    "D49" = "56.08"
    "E49" = "  +2.99%  "
    "E50" = "  +5.05%  "
    "D51" = "8.11"
    "E51" = "  -0.32%  "
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
